$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain plain text so numeric-looking strings
# (e.g. "0.999") are not auto-converted to Number cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '49.777.90'
$ws.Range('E2').Value = '  +3.30%  '
$ws.Range('D3').Value = '2.616.87'
$ws.Range('E3').Value = '  +4.58%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').Value = '110.28'
$ws.Range('E5').Value = '  +2.03%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '324.06'
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('D7').Value = '0.535'
$ws.Range('E7').Value = '  +1.73%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = '0.561'
$ws.Range('E9').Value = '  +3.79%  '
$ws.Range('D10').Value = '40.89'
$ws.Range('E10').Value = '  +2.53%  '
$ws.Range('D11').Value = '20.67'
$ws.Range('E11').Value = '  +2.34%  '
$ws.Range('D12').Value = '0.0824'
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('E14').Value = '  +2.18%  '
$ws.Range('D15').Value = '3.028.68'
$ws.Range('E15').Value = '  +4.69%  '
$ws.Range('D16').Value = '2.610.58'
$ws.Range('E16').Value = '  +4.25%  '
$ws.Range('D17').Value = '0.874'
$ws.Range('E17').Value = '  +3.49%  '
$ws.Range('D18').Value = '49.679.19'
$ws.Range('E18').Value = '  +3.44%  '
$ws.Range('E19').Value = '  +11.74%  '
$ws.Range('D20').Value = '13.35'
$ws.Range('E20').Value = '  +2.03%  '
$ws.Range('D21').Value = '6.78'
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('D22').Value = '0.0₃0955'
$ws.Range('E22').Value = '  +0.94%  '
$ws.Range('D23').Value = '281.52'
$ws.Range('E23').Value = '  +1.44%  '
$ws.Range('D24').Value = '72.82'
$ws.Range('E24').Value = '  +1.18%  '
$ws.Range('D25').Value = '2.58'
$ws.Range('E25').Value = '  +0.96%  '
$ws.Range('D26').Value = '26.65'
$ws.Range('E26').Value = '  +3.49%  '
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('E28').Value = '  -6.92%  '
$ws.Range('D29').Value = '9.99'
$ws.Range('E29').Value = '  +2.16%  '
$ws.Range('E30').Value = '  +3.32%  '
$ws.Range('D31').Value = '36.35'
$ws.Range('E31').Value = '  +3.14%  '
$ws.Range('D32').Value = '49.61'
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('D33').Value = '5.47'
$ws.Range('E33').Value = '  +2.43%  '
$ws.Range('D34').Value = '19.64'
$ws.Range('E34').Value = '  +0.47%  '
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('D36').Value = '0.0797'
$ws.Range('E36').Value = '  +1.55%  '
$ws.Range('D37').Value = '2.06'
$ws.Range('E37').Value = '  +5.23%  '
$ws.Range('E38').Value = '  +2.88%  '
$ws.Range('D39').Value = '3.09'
$ws.Range('E39').Value = '  +6.07%  '
$ws.Range('D40').Value = '22.93'
$ws.Range('E40').Value = '  +6.81%  '
$ws.Range('D41').Value = '0.113'
$ws.Range('E41').Value = '  +1.04%  '
$ws.Range('D42').Value = '123.27'
$ws.Range('E42').Value = '  +1.78%  '
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('D44').Value = '0.0316'
$ws.Range('E44').Value = '  +4.04%  '
$ws.Range('D45').Value = '3.37'
$ws.Range('E45').Value = '  +6.80%  '
$ws.Range('D46').Value = '2.051.69'
$ws.Range('E46').Value = '  +2.31%  '
$ws.Range('D47').Value = '2.22'
$ws.Range('E47').Value = '  +11.60%  '
$ws.Range('D48').Value = '2.02'
$ws.Range('E48').Value = '  +9.25%  '
$ws.Range('D49').Value = '9.02'
$ws.Range('E49').Value = '  +0.75%  '
$ws.Range('D50').Value = '5.38'
$ws.Range('E50').Value = '  +3.95%  '
$ws.Range('D51').Value = '82.19'
$ws.Range('E51').Value = '  +2.76%  '

# Drop the temporary text-number-format styling so the D column
# cells end up with no explicit style, matching the source data.
$ws.Range("D2:D51").ClearFormats()
